# "Update countries & provincias Spain"
# Applies the 1-Jul-2020 02:22 refresh of the COVID country table:
#  - reorders two pairs of country names (Libia/Suazilandia and Laos/Santa
#    Lucia + Fiyi/Dominica swap back into alphabetic-ish ranking order)
#  - refreshes the numeric columns (Casos totales, Nuevos casos, Casos
#    activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for the
#    rows whose figures moved
#  - bumps the "Datos actualizados" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder country names (A column) ---------------------------------
# Andorra(142), Suazilandia(143), Libia(144), Liberia(145)
#   -> Andorra(142), Libia(143), Suazilandia(144), Liberia(145)
$ws.Range("A143").Value = "Libia"
$ws.Range("A144").Value = "Suazilandia"

# Nueva Caledonia(202), Laos(203), Santa Lucia(204), Fiyi(205),
# Dominica(206), Islas Virgenes de los Estados Unidos(207)
#   -> Nueva Caledonia(202), Santa Lucia(203), Laos(204), Dominica(205),
#      Fiyi(206), Islas Virgenes de los Estados Unidos(207)
$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("A204").Value = "Laos"
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"

# --- Refresh timestamp string -------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Julio de 2020 a las 02:22"

# --- Refresh numeric data -------------------------------------------------
# row 4: Estados Unidos
$ws.Range("B4").Value = 2727305
$ws.Range("C4").Value = 45494
$ws.Range("D4").Value = 1140890
$ws.Range("E4").Value = 1456304
$ws.Range("G4").Value = 753
$ws.Range("H4").Value = 130111

# row 5: Brasil
$ws.Range("B5").Value = 1408485
$ws.Range("C5").Value = 37997
$ws.Range("E5").Value = 558789
$ws.Range("G5").Value = 1271
$ws.Range("H5").Value = 59656

# row 17: Alemania
$ws.Range("B17").Value = 195832
$ws.Range("C17").Value = 440
$ws.Range("E17").Value = 7680

# row 22: Canada
$ws.Range("B22").Value = 104204
$ws.Range("C22").Value = 286
$ws.Range("D22").Value = 67594
$ws.Range("E22").Value = 28019

# row 43: Panama
$ws.Range("B43").Value = 33550
$ws.Range("C43").Value = 765
$ws.Range("D43").Value = 15745
$ws.Range("E43").Value = 17174
$ws.Range("G43").Value = 11
$ws.Range("H43").Value = 631

# row 133: Niger
$ws.Range("D133").Value = 943
$ws.Range("E133").Value = 65

# row 137: Uruguay
$ws.Range("B137").Value = 936
$ws.Range("C137").Value = 4
$ws.Range("D137").Value = 824
$ws.Range("E137").Value = 85

# row 143: now Libia (after the swap above)
$ws.Range("B143").Value = 824
$ws.Range("C143").Value = 22
$ws.Range("D143").Value = 209
$ws.Range("E143").Value = 591
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 24

# row 144: now Suazilandia (after the swap above)
$ws.Range("B144").Value = 812
$ws.Range("C144").Value = 17
$ws.Range("D144").Value = 408
$ws.Range("E144").Value = 393
$ws.Range("H144").Value = 11

# row 185: Seychelles
$ws.Range("B185").Value = 81
$ws.Range("C185").Value = 4
$ws.Range("E185").Value = 70
